$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new class D22CQCN01-N, year 2022-2023, advisor Phạm Bình An
$ws.Range("A6").Value = 1026
$ws.Range("B6").Value = "D22CQCN01-N"
$ws.Range("C6").Value = "2022-2023"
$ws.Range("D6").Value = "Phạm Bình An"
$ws.Range("E6").Value = 1

# Row 7: same class D22CQCN01-N, year 2023-2024, advisor Đinh Anh
$ws.Range("A7").Value = 1027
$ws.Range("B7").Value = "D22CQCN01-N"
$ws.Range("C7").Value = "2023-2024"
$ws.Range("D7").Value = "Đinh Anh"
$ws.Range("E7").Value = 1
